$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Fix typo in row 12's activity text (add missing comma)
# ------------------------------------------------------------------
$ws.Cells.Item(12, 2).Value2 = "Discussion, planification et organisation de groupe"

# ------------------------------------------------------------------
# 2. Re-arrange the journal entries around rows 15-17:
#    - old row 15 ("Ajout de l'action...", 43186, 3h, tall row) moves to row 16
#      with its date corrected to 43185
#    - row 15 becomes a brand-new entry ("Discussion, planification...",
#      43185, 1.5h) at normal row height
#    - row 17 (previously empty) becomes a brand-new entry
#      ("Ajout du code...", 43186, 1.5h)
# ------------------------------------------------------------------

# Move old row 15 content down into row 16 (copy first, then fix the date)
$ws.Cells.Item(16, 1).Value2 = $ws.Cells.Item(15, 1).Value2
$ws.Cells.Item(16, 2).Value2 = $ws.Cells.Item(15, 2).Value2
$ws.Cells.Item(16, 3).Value2 = $ws.Cells.Item(15, 3).Value2
$ws.Rows(16).RowHeight = $ws.Rows(15).RowHeight
$ws.Cells.Item(16, 1).Value2 = 43185

# New row 15 content
$ws.Cells.Item(15, 1).Value2 = 43185
$ws.Cells.Item(15, 2).Value2 = "Discussion, planification et organisation de groupe"
$ws.Cells.Item(15, 3).Value2 = 1.5
$ws.Rows(15).AutoFit()

# New row 17 content
$ws.Cells.Item(17, 1).Value2 = 43186
$ws.Cells.Item(17, 2).Value2 = "Ajout du code gérant la sauvegarde et l'ouverture de fichier"
$ws.Cells.Item(17, 3).Value2 = 1.5

# ------------------------------------------------------------------
# 3. The journal now has one more filled-in row, so the blank buffer
#    rows shrink by one: what used to be the "Total" row (32) becomes
#    a blank row, and "Total" moves down to row 33, with the SUM
#    formula extended to include the new row.
# ------------------------------------------------------------------

# Copy row 32 (currently the Total row) down into row 33, then overwrite
# with the updated label/formula (range now C5:C32).
$ws.Range("A32:C32").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)
$ws.Cells.Item(33, 2).Value2 = "Total"
$ws.Cells.Item(33, 3).Formula = "=SUM(C5:C32)"

# Turn row 32 back into a blank buffer row matching the others.
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)
$ws.Range("A32:C32").ClearContents()
$ws.Cells.Item(33, 1).Clear()

# ------------------------------------------------------------------
# 4. Update the current selection to match the author's last position.
# ------------------------------------------------------------------
$ws.Range("C18").Select()
